$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RedMediumBorder($range, $edge) {
    $b = $range.Borders.Item($edge)
    $b.Color = 255
    $b.LineStyle = 1
    $b.Weight = -4138
}

# ---------------------------------------------------------------------------
# 1. New columns G (Tuan 6) and H (Tuan 7) appended to the existing summary
#    table (rows 3-6).
# ---------------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 17.36
$ws.Columns("H").ColumnWidth = 17.93

# Header row - copy the green header format from A3 onto G3:H3
$ws.Range("A3").Copy() | Out-Null
$ws.Range("G3:H3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G3").Value = "Tuần 6"
$ws.Range("H3").Value = "Tuần 7"

$ws.Range("G4").Value = "1. Chức năng đăng tin (video) để promotion về địa điểm du lịch của Tour Guide"
Set-RedMediumBorder $ws.Range("G4") 7    # left

$ws.Range("H4").Value = "2. Chức năng Chat.`n3. Chức năng Book Tour."
Set-RedMediumBorder $ws.Range("H4") 7    # left
$ws.Range("H4").WrapText = $true

$ws.Range("G5").Value = "2. Quản lý tin đăng (Tour Guide) (Manage seft tour).`n"
$ws.Range("G5").WrapText = $true

$ws.Range("H5").Value = "1. Chức năng View Tour + Search Tour.`n4. Chức năng FeedBack."
Set-RedMediumBorder $ws.Range("H5") 7    # left
$ws.Range("H5").WrapText = $true

# F6 fill: yellow -> theme background 1 (white)
$ws.Range("F6").Interior.ThemeColor = 2

# ---------------------------------------------------------------------------
# 2. Spacer rows
# ---------------------------------------------------------------------------
$ws.Rows(9).RowHeight = 15
$ws.Rows(14).RowHeight = 15

# ---------------------------------------------------------------------------
# 3. "Tuan 6" box: A10:D12
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Tuần 6:"
$ws.Range("A11").Value = "1. Chức năng đăng tin (video) để promotion về địa điểm du lịch của Tour Guide"
$ws.Range("A12").Value = "2. Quản lý tin đăng (Tour Guide) (Manage seft tour)"

Set-RedMediumBorder $ws.Range("A10:D10") 8   # top
Set-RedMediumBorder $ws.Range("A10:A12") 7   # left
Set-RedMediumBorder $ws.Range("D10:D12") 10  # right
Set-RedMediumBorder $ws.Range("A12:D12") 9   # bottom

# ---------------------------------------------------------------------------
# 4. "Tuan 7" box: A15:B19
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Tuần 7 :"
$ws.Range("A16").Value = "1. Chức năng View Tour + Search Tour."
$ws.Range("A17").Value = "2. Chức năng Chat"
$ws.Range("A18").Value = "3. Chức năng Book Tour"
$ws.Range("A19").Value = "4. Chức năng FeedBack"

Set-RedMediumBorder $ws.Range("A15:B15") 8   # top
Set-RedMediumBorder $ws.Range("A15:A19") 7   # left
Set-RedMediumBorder $ws.Range("B15:B19") 10  # right
Set-RedMediumBorder $ws.Range("A19:B19") 9   # bottom

# ---------------------------------------------------------------------------
# 5. "Tuan 8" plain block
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Tuần 8 :"
$ws.Range("A22").Value = "1. QL User"
$ws.Range("B22").Value = "Admin"
$ws.Range("A23").Value = "2. QL Tour"

# ---------------------------------------------------------------------------
# 6. Selection, matching the authored edit
# ---------------------------------------------------------------------------
$ws.Range("G4").Select() | Out-Null

Write-Host "edit applied"
